# Apply corrected vaccination/attack numbers to the measles_costs sheet,
# then touch the selection so it matches the reviewed workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("measles_costs")

# Corrected "total cases" (col B) and "unvaccinated cases" (col D) figures
# for each district row (19-38). Columns C (vaccinated) and E (GP visits)
# are unchanged but are re-written too so every cell in the block picks up
# the same refreshed formatting as the cells whose numbers actually moved.
$rows = @(
    @{ Row = 19; B = 436347; C = 47930; D = 24719; E = 13840 },
    @{ Row = 20; B = 205998; C = 19642; D = 6641;  E = 3548  },
    @{ Row = 21; B = 482178; C = 48190; D = 19384; E = 10520 },
    @{ Row = 22; B = 283704; C = 30414; D = 14853; E = 8250  },
    @{ Row = 23; B = 469299; C = 52672; D = 28389; E = 16008 },
    @{ Row = 24; B = 151695; C = 14842; D = 5550;  E = 2991  },
    @{ Row = 25; B = 138375; C = 14288; D = 6339;  E = 3477  },
    @{ Row = 26; B = 98196;  C = 10015; D = 4289;  E = 2343  },
    @{ Row = 27; B = 162561; C = 16432; D = 6851;  E = 3732  },
    @{ Row = 28; B = 136995; C = 12320; D = 3086;  E = 1617  },
    @{ Row = 29; B = 151686; C = 14257; D = 4526;  E = 2407  },
    @{ Row = 30; B = 55620;  C = 4949;  D = 1156;  E = 604   },
    @{ Row = 31; B = 297420; C = 29797; D = 12080; E = 6561  },
    @{ Row = 32; B = 43650;  C = 4557;  D = 2083;  E = 1147  },
    @{ Row = 33; B = 109752; C = 10821; D = 4158;  E = 2247  },
    @{ Row = 34; B = 359313; C = 37303; D = 16798; E = 9232  },
    @{ Row = 35; B = 41112;  C = 3744;  D = 1011;  E = 532   },
    @{ Row = 36; B = 525549; C = 54788; D = 24948; E = 13729 },
    @{ Row = 37; B = 32151;  C = 3009;  D = 936;   E = 497   },
    @{ Row = 38; B = 60120;  C = 5772;  D = 2008;  E = 1075  }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("B$n").Value = $r.B
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
}

# Recalculate so every dependent formula (G:M, N:W, row 40 totals, etc.)
# picks up the corrected figures.
$excel.CalculateFull()

# The reviewer's cursor ended up on G42 instead of P42.
$ws.Range("G42").Select()
